$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 123
$ws.Range("H123").Value = 69780
$ws.Range("J123").Value = 69780
$ws.Range("L123").Value = 69780
$ws.Range("N123").Value = -79580
# Row 129
$ws.Range("H129").Value = 1136.1578
$ws.Range("J129").Value = 1235.4412
$ws.Range("L129").Value = 3706.3236
$ws.Range("N129").Value = -13706.3236
# Row 137
$ws.Range("H137").Value = 20409374
$ws.Range("I137").Value = 30303942
$ws.Range("J137").Value = 1824.125
$ws.Range("K137").Value = 90911826
$ws.Range("L137").Value = 5472.375
$ws.Range("M137").Value = -90909276
$ws.Range("N137").Value = -10572.375
# Row 138
$ws.Range("H138").Value = 9224182
$ws.Range("I138").Value = 2234723.8
$ws.Range("J138").Value = 15154631
$ws.Range("K138").Value = 6704171.399999999
$ws.Range("L138").Value = 45463893
$ws.Range("M138").Value = -6699031.399999999
$ws.Range("N138").Value = -45474173
# Row 141
$ws.Range("H141").Value = 4044.7021
$ws.Range("I141").Value = 2269.2058
$ws.Range("J141").Value = 8688.308000000001
$ws.Range("K141").Value = 6807.617400000001
$ws.Range("L141").Value = 26064.924
$ws.Range("M141").Value = -1627.617400000001
$ws.Range("N141").Value = -36424.924

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2033.8572
$ws.Range("I61").Value = 1301.0714
$ws.Range("K61").Value = 1301.0714
$ws.Range("M61").Value = -1089.0714
# Row 122
$ws.Range("H122").Value = 9580.643
$ws.Range("I122").Value = 10635.75
$ws.Range("J122").Value = 3250
$ws.Range("K122").Value = 31907.25
$ws.Range("L122").Value = 9750
$ws.Range("M122").Value = -29457.25
$ws.Range("N122").Value = -14650
# Row 132
$ws.Range("H132").Value = 2014.95
$ws.Range("I132").Value = 1569.4509
$ws.Range("J132").Value = 4539.4443
$ws.Range("K132").Value = 4708.3527
$ws.Range("L132").Value = 13618.3329
$ws.Range("M132").Value = -2178.3527
$ws.Range("N132").Value = -18678.3329
# Row 136
$ws.Range("H136").Value = 2033.8572
$ws.Range("I136").Value = 1301.0714
$ws.Range("K136").Value = 3903.2142
$ws.Range("M136").Value = -1353.2142

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 3401.2
$ws.Range("I105").Value = 3223.4783
$ws.Range("J105").Value = 3741.8333
$ws.Range("K105").Value = 3223.4783
$ws.Range("L105").Value = 3741.8333
$ws.Range("M105").Value = -1476.4783
$ws.Range("N105").Value = -7235.8333
# Row 134
$ws.Range("H134").Value = 2458.1538
$ws.Range("I134").Value = 1647.35
$ws.Range("J134").Value = 5160.8335
$ws.Range("K134").Value = 4942.049999999999
$ws.Range("L134").Value = 15482.5005
$ws.Range("M134").Value = -2407.049999999999
$ws.Range("N134").Value = -20552.5005

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 47
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
# Row 58
$ws.Range("H58").Value = 1642.2069
$ws.Range("I58").Value = 1119.28
$ws.Range("K58").Value = 1119.28
$ws.Range("M58").Value = -916.28
# Row 86
$ws.Range("H86").Value = 125003500
$ws.Range("J86").Value = 4666.6665
$ws.Range("L86").Value = 4666.6665
$ws.Range("N86").Value = -6912.6665
# Row 89
$ws.Range("H89").Value = 125003500
$ws.Range("J89").Value = 4666.6665
$ws.Range("L89").Value = 23333.3325
$ws.Range("N89").Value = -34565.3325
# Row 132
$ws.Range("H132").Value = 1415.4565
$ws.Range("I132").Value = 1025.8049
$ws.Range("J132").Value = 4610.6
$ws.Range("K132").Value = 3077.4147
$ws.Range("L132").Value = 13831.8
$ws.Range("M132").Value = -547.4147000000003
$ws.Range("N132").Value = -18891.8
# Row 134
$ws.Range("H134").Value = 2082.0784
$ws.Range("I134").Value = 1295.159
$ws.Range("J134").Value = 7028.4287
$ws.Range("K134").Value = 3885.477
$ws.Range("L134").Value = 21085.2861
$ws.Range("M134").Value = -1350.477
$ws.Range("N134").Value = -26155.2861
# Row 136
$ws.Range("H136").Value = 1642.2069
$ws.Range("I136").Value = 1119.28
$ws.Range("K136").Value = 3357.84
$ws.Range("M136").Value = -807.8400000000001

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 5051848
$ws.Range("I131").Value = 367.8
$ws.Range("J131").Value = 5953898
$ws.Range("K131").Value = 1103.4
$ws.Range("L131").Value = 17861694
$ws.Range("M131").Value = 3936.6
$ws.Range("N131").Value = -17871774

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5315.075
$ws.Range("I70").Value = 5279.788
$ws.Range("J70").Value = 5481.4287
$ws.Range("K70").Value = 5279.788
$ws.Range("L70").Value = 5481.4287
$ws.Range("M70").Value = -5009.788
$ws.Range("N70").Value = -6021.4287
# Row 73
$ws.Range("H73").Value = 5315.075
$ws.Range("I73").Value = 5279.788
$ws.Range("J73").Value = 5481.4287
$ws.Range("K73").Value = 5279.788
$ws.Range("L73").Value = 5481.4287
$ws.Range("M73").Value = -4343.788
$ws.Range("N73").Value = -7353.4287
# Row 80
$ws.Range("H80").Value = 2607.1428
$ws.Range("I80").Value = 2512.5
$ws.Range("J80").Value = 2910
$ws.Range("K80").Value = 2512.5
$ws.Range("L80").Value = 2910
$ws.Range("M80").Value = -1514.5
$ws.Range("N80").Value = -4906
# Row 83
$ws.Range("H83").Value = 2607.1428
$ws.Range("I83").Value = 2512.5
$ws.Range("J83").Value = 2910
$ws.Range("K83").Value = 12562.5
$ws.Range("L83").Value = 14550
$ws.Range("M83").Value = -7570.5
$ws.Range("N83").Value = -24534
# Row 122
$ws.Range("H122").Value = 1284.6666
$ws.Range("I122").Value = 570.8
$ws.Range("J122").Value = 2177
$ws.Range("K122").Value = 1712.4
$ws.Range("L122").Value = 6531
$ws.Range("M122").Value = 737.6000000000001
$ws.Range("N122").Value = -11431

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
# Row 27
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()
# Row 136
$ws.Range("H136").Value = 7154.846
$ws.Range("I136").Value = 2995.4443
$ws.Range("J136").Value = 16513.5
$ws.Range("K136").Value = 8986.332900000001
$ws.Range("L136").Value = 49540.5
$ws.Range("M136").Value = -6436.332900000001
$ws.Range("N136").Value = -54640.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 16892.309
$ws.Range("I62").Value = 16544.445
$ws.Range("J62").Value = 17675
$ws.Range("K62").Value = 16544.445
$ws.Range("L62").Value = 17675
$ws.Range("M62").Value = -15920.445
$ws.Range("N62").Value = -18923
# Row 65
$ws.Range("H65").Value = 16892.309
$ws.Range("I65").Value = 16544.445
$ws.Range("J65").Value = 17675
$ws.Range("K65").Value = 82722.22500000001
$ws.Range("L65").Value = 88375
$ws.Range("M65").Value = -79602.22500000001
$ws.Range("N65").Value = -94615
# Row 100
$ws.Range("H100").Value = 1245.5
$ws.Range("I100").Value = 1160.6666
$ws.Range("J100").Value = 1500
$ws.Range("K100").Value = 2321.3332
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1780.3332
$ws.Range("N100").Value = -4082
# Row 122
$ws.Range("H122").Value = 1431.7333
$ws.Range("I122").Value = 1152.5
$ws.Range("J122").Value = 1990.2
$ws.Range("K122").Value = 3457.5
$ws.Range("L122").Value = 5970.6
$ws.Range("M122").Value = -1007.5
$ws.Range("N122").Value = -10870.6
# Row 123
$ws.Range("H123").Value = 32966.668
$ws.Range("J123").Value = 32966.668
$ws.Range("L123").Value = 32966.668
$ws.Range("N123").Value = -42766.668
# Row 126
$ws.Range("H126").Value = 103100.4
$ws.Range("I126").Value = 114244.89
$ws.Range("J126").Value = 2800
$ws.Range("K126").Value = 342734.67
$ws.Range("L126").Value = 8400
$ws.Range("M126").Value = -340264.67
$ws.Range("N126").Value = -13340
# Row 132
$ws.Range("H132").Value = 3145.6191
$ws.Range("I132").Value = 3002.7188
$ws.Range("J132").Value = 3602.9
$ws.Range("K132").Value = 9008.1564
$ws.Range("L132").Value = 10808.7
$ws.Range("M132").Value = -6478.1564
$ws.Range("N132").Value = -15868.7
